$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reverse the "Periodo Mora" period values for rows 16-22 (E column)
$ws.Range("E16").Value = "2411"
$ws.Range("E17").Value = "2410"
$ws.Range("E18").Value = "2409"
$ws.Range("E19").Value = "2408"
$ws.Range("E20").Value = "2407"
$ws.Range("E21").Value = "2406"
$ws.Range("E22").Value = "2405"

# Swap the "Valor Mora" amounts between the first and last rows
$ws.Range("F16").Value = 24266
$ws.Range("F22").Value = 52000
